$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I: new "Other found locations" header + values ---
$ws.Range("I1").Value = "Other found locations"
$ws.Range("I2").Value = ""
$ws.Range("I3").Value = "_PMC_Springer"
$ws.Range("I4").Value = "_PMC_Springer"
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("I10").Value = "_PMC"
$ws.Range("I11").Value = ""
$ws.Range("I12").Value = "_PMC_elsevier"
$ws.Range("I13").Value = ""
$ws.Range("I14").Value = ""
$ws.Range("I15").Value = "_PMC_elsevier"
$ws.Range("I16").Value = ""
$ws.Range("I17").Value = "_PMC"
$ws.Range("I18").Value = "_PMC"

# --- Column E: updated Authors text ---
$ws.Range("E3").Value = "[Ana Rosa%Linde-Arias%linde14@yahoo.com%1,   Maria%Roura%NULL%1,   Eduardo%Siqueira%NULL%1]"
$ws.Range("E4").Value = "[Maria%Belizan%mbelizan@iecs.org.ar%1,   Edna%Maradiaga%edjamar3006@yahoo.com%1,   Javier%Roberti%jroberti@iecs.org.ar%1,   Maricela%Casco-Aguilar%marykasco@yahoo.com%1,   Alison F.%Ortez%alison_fabiola@yahoo.es%1,   Juan C.%Avila-Flores%javilaflores3@gmail.com%1,   Gloria%González%marilyntoin@yahoo.com%1,   Carolina%Bustillo%mcbu1502@yahoo.com%1,   Alejandra%Calderón%lilianalecalderon@gmail.com%1,   Harry%Bock%hbockme@hotmail.com%1,   María L.%Cafferata%NULL%1,   Adriano B.%Tavares%adriano.b.tavares@gmail.com%1,   Jackeline%Alger%jackelinealger@gmail.com%1,   Moazzam%Ali%alimoa@who.int%1]"
$ws.Range("E5").Value = "[ Catherine A.%Boyd%null%2,    Julie A.%Gazmararian%null%1,    Winifred Wilkins%Thompson%null%1,  Catherine A.%Boyd%null%0,  Julie A.%Gazmararian%null%1,  Winifred Wilkins%Thompson%null%1]"
$ws.Range("E6").Value = "[Benner%P.%coreGivesNoEmail%1,  Bohme%G.%coreGivesNoEmail%1,  Chiang%H-H%coreGivesNoEmail%4,  Chiang%H-H%coreGivesNoEmail%0,  Chiang%H-H%coreGivesNoEmail%0,  Chiang%H-H%coreGivesNoEmail%0,  Donohoe%J.%coreGivesNoEmail%1,  Hsien-Hsien%Chiang%coreGivesNoEmail%1,  I-Ling%Sue%coreGivesNoEmail%1,  International%Council of Nurses%coreGivesNoEmail%1,  Mei-Bih%Chen%coreGivesNoEmail%1,  Morgan%DL%coreGivesNoEmail%1,  Packer%MJ%coreGivesNoEmail%1]"
$ws.Range("E7").Value = "[Chee%YO%coreGivesNoEmail%1,  Dodgson%JE%coreGivesNoEmail%1,  Tarrant%M%coreGivesNoEmail%1,  Watkins%A%coreGivesNoEmail%1]"
$ws.Range("E8").Value = "[ Elin%Erland%null%2,    Bente%Dahl%null%1,  Elin%Erland%null%0,  Bente%Dahl%null%1]"
$ws.Range("E9").Value = "[Davis%Mark%coreGivesNoEmail%1,  Flowers%Paul%coreGivesNoEmail%1,  Lohm%Davina%coreGivesNoEmail%1,  Stephenson%Niamh%coreGivesNoEmail%1,  Waller%Emily%coreGivesNoEmail%1]"
$ws.Range("E10").Value = "[Hector M.%Gomez%NULL%1,   Carlos%Mejia Arbelaez%NULL%1,   Jovana A.%Ocampo Cañas%ja.ocampo@uniandes.edu.co%1]"
$ws.Range("E11").Value = "[ Susan%Jones%null%2,    Betty%Sam%null%1,    Florence%Bull%null%1,    Steven Bagie%Pieh%null%1,    Jaki%Lambert%null%1,    Florence%Mgawadere%null%1,    Somasundari%Gopalakrishnan%null%1,    Charles A.%Ameh%null%1,    Nynke%van den Broek%null%1,  Susan%Jones%null%0,  Betty%Sam%null%1,  Florence%Bull%null%1,  Steven Bagie%Pieh%null%1,  Jaki%Lambert%null%1,  Florence%Mgawadere%null%1,  Somasundari%Gopalakrishnan%null%1,  Charles A.%Ameh%null%1,  Nynke%van den Broek%null%1]"
$ws.Range("E12").Value = "[Caroline S.E.%Homer%NULL%1,   Miranda%Davies-Tuck%NULL%1,   Hannah G.%Dahlen%NULL%1,   Vanessa L.%Scarf%NULL%1]"
$ws.Range("E14").Value = "[ Molly M.%Lynch%null%2,    Elizabeth W.%Mitchell%null%1,    Jennifer L.%Williams%null%1,    Kelly%Brumbaugh%null%1,    Michelle%Jones-Bell%null%1,    Debra E.%Pinkney%null%1,    Christine M.%Layton%null%1,    Patricia W.%Mersereau%null%1,    Juliette S.%Kendrick%null%1,    Paula Eguino%Medina%null%1,    Lucia Rojas%Smith%null%1,  Molly M.%Lynch%null%0,  Elizabeth W.%Mitchell%null%1,  Jennifer L.%Williams%null%1,  Kelly%Brumbaugh%null%1,  Michelle%Jones-Bell%null%1,  Debra E.%Pinkney%null%1,  Christine M.%Layton%null%1,  Patricia W.%Mersereau%null%1,  Juliette S.%Kendrick%null%1,  Paula Eguino%Medina%null%1,  Lucia Rojas%Smith%null%1]"
$ws.Range("E15").Value = "[Niuniu%Sun%NULL%1,   Luoqun%Wei%NULL%1,   Suling%Shi%NULL%1,   Dandan%Jiao%NULL%1,   Runluo%Song%NULL%1,   Lili%Ma%NULL%1,   Hongwei%Wang%NULL%2,   Chao%Wang%NULL%1,   Zhaoguo%Wang%NULL%1,   Yanli%You%NULL%1,   Shuhua%Liu%NULL%1,   Hongyun%Wang%NULL%1]"
$ws.Range("E17").Value = "[Veronika%Tirado%Veronika.Tirado@ki.se%1,   Santiago A.%Morales Mesa%NULL%1,   John%Kinsman%NULL%1,   Anna Mia%Ekström%NULL%0,   Berta N.%Restrepo Jaramillo%NULL%1]"
$ws.Range("E18").Value = "[Simon N%Williams%NULL%1,   Christopher J%Armitage%NULL%2,   Christopher J%Armitage%NULL%0,   Tova%Tampe%NULL%1,   Kimberly%Dienes%NULL%1]"

# --- Column F/G: ID / ID Format updates for rows 6,7,9 ---
$ws.Range("F6").Value = "not found"
$ws.Range("F7").Value = "not found"
$ws.Range("F9").Value = "not found"
$ws.Range("G6").Value = "N/A"
$ws.Range("G7").Value = "N/A"
$ws.Range("G9").Value = "N/A"
